$p = $ppt.ActivePresentation

# Slide 2: TextBox 3 shape has "The" + " " + "Moon" split across 3 runs;
# merge into a single run "The Moon". Setting directly to the same
# concatenated text is a no-op for the underlying run structure, so we
# first set a distinct placeholder value to force the merge, then set
# the real text.
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(2)
$tb2.TextFrame.TextRange.Text = "__tmp__"
$tb2.TextFrame.TextRange.Text = "The Moon"

# Slide 3: Title 1 shape has "One" + " " + "More" split across 3 runs;
# merge into a single run "One More".
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "__tmp__"
$title3.TextFrame.TextRange.Text = "One More"

# Slide 3: TextBox 3 shape has "The" + " " + "Moon" split across 3 runs;
# merge into a single run "The Moon".
$tb3 = $s3.Shapes.Item(3)
$tb3.TextFrame.TextRange.Text = "__tmp__"
$tb3.TextFrame.TextRange.Text = "The Moon"
